$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at A (old A -> B, old B -> C); formatting of the old
#    columns rides along automatically.
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()

# ---------------------------------------------------------------------------
# 2. Give the brand-new column A the same base formatting (font/border) as
#    its row neighbour before we touch the fill colours, so we never mint a
#    duplicate font entry - only fills/alignment differ going forward.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A5").PasteSpecial(-4122) | Out-Null

$ws.Range("B3").Copy() | Out-Null
$ws.Range("A6:A10").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Column A values ("Acteurs").
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Acteurs"
$ws.Range("A2").Value = "SGP"
$ws.Range("A3").Value = "SGC"
$ws.Range("A4").Value = "Personne morale"
$ws.Range("A5").Value = "Personne physique"
$ws.Range("A6").Value = "Client non enregistré"
$ws.Range("A7").Value = "Client non enregistré"
$ws.Range("A8").Value = "SGC"
$ws.Range("A9").Value = "SGC"
$ws.Range("A10").Value = "SGP"

# ---------------------------------------------------------------------------
# 4. Re-colour the header row (A1:C1) with the new themed dark fill.
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").Interior.ThemeColor = 3

# ---------------------------------------------------------------------------
# 5. Re-colour the "category" cells (col A + col B) that used to be plain
#    light-blue so they pick up the new themed blues. Rows that were yellow
#    stay yellow (same colour, new fill slot - no action required).
# ---------------------------------------------------------------------------
$ws.Range("A3").Interior.ThemeColor = 5
$ws.Range("A6:A10").Interior.ThemeColor = 5

$ws.Range("B3").Interior.ThemeColor = 5
$ws.Range("B6:B10").Interior.ThemeColor = 5

# ---------------------------------------------------------------------------
# 6. Column C (the old "Commentaires" column) : left-align everything and
#    recolour to match its row (yellow rows stay yellow, others go themed
#    medium blue).
# ---------------------------------------------------------------------------
$ws.Range("C2:C10").HorizontalAlignment = -4131

$ws.Range("C2").Interior.Color = 65535
$ws.Range("C4").Interior.Color = 65535
$ws.Range("C5").Interior.Color = 65535

$ws.Range("C3").Interior.ThemeColor = 5
$ws.Range("C6:C10").Interior.ThemeColor = 5

# ---------------------------------------------------------------------------
# 7. Column widths / dimensions.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 19.44140625
$ws.Columns("B").ColumnWidth = 45.44140625
$ws.Columns("C").ColumnWidth = 61.88671875

# ---------------------------------------------------------------------------
# 8. Misc view state.
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()
